# B1--and-B2-PowerPoint.pptx edit
#
# The table on slide 5 ("B1- TYPES OF FINANCIAL DOCUMENTS") has its table
# style (Table Design gallery selection) changed from the document's
# custom "Table_0" style to a different (built-in) table style.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)

if ($sh.HasTable) {
    $sh.Table.ApplyStyle("{62529EB1-D295-4DAA-A69C-7FD25FFCD3A7}")
}
